# Update Betfair Back/Lay odds figures for 2026-01-14 (data refresh).
# Each statement below writes the single updated odds/price value for one cell,
# matching the refreshed snapshot of the "Jogos do Dia" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4.8
$ws.Range("G2").Value = 5.1
$ws.Range("H2").Value = 1.86
$ws.Range("I2").Value = 1.89
$ws.Range("J2").Value = 3.7
$ws.Range("K2").Value = 3.8
$ws.Range("L2").Value = 1.48
$ws.Range("M2").Value = 1.09
$ws.Range("N2").Value = 3.3
$ws.Range("O2").Value = 1.41
$ws.Range("P2").Value = 1.78
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 4.1
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 1.9
$ws.Range("V2").Value = 2.1
$ws.Range("W2").Value = 1.25
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 10.5
$ws.Range("AA2").Value = 20
$ws.Range("AB2").Value = 15
$ws.Range("AC2").Value = 8.4
$ws.Range("AE2").Value = 22
$ws.Range("AF2").Value = 36
$ws.Range("AG2").Value = 20
$ws.Range("AH2").Value = 22
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 130
$ws.Range("AL2").Value = 85
$ws.Range("AM2").Value = 140
$ws.Range("AN2").Value = 110
$ws.Range("AO2").Value = 15.5

$ws.Range("X5").Value = 970
$ws.Range("Y5").Value = 30
$ws.Range("AB5").Value = 9.800000000000001
$ws.Range("AD5").Value = 36
$ws.Range("AE5").Value = 140
$ws.Range("AF5").Value = 10.5
$ws.Range("AI5").Value = 120
$ws.Range("AJ5").Value = 14.5
$ws.Range("AK5").Value = 19
$ws.Range("AO5").Value = 180

$ws.Range("F6").Value = 2.1
$ws.Range("I6").Value = 22
$ws.Range("J6").Value = 1.03
$ws.Range("K6").Value = 25
$ws.Range("N6").Value = 2.3
$ws.Range("O6").Value = 1.1
$ws.Range("P6").Value = 2.42
$ws.Range("R6").Value = 2.3
$ws.Range("S6").Value = 1.69
$ws.Range("V6").Value = 1.04

$ws.Range("F7").Value = 1.99
$ws.Range("K7").Value = 3.85

$ws.Range("F8").Value = 1.34
$ws.Range("I8").Value = 13.5
$ws.Range("AD8").Value = 50
$ws.Range("AH8").Value = 40
$ws.Range("AK8").Value = 16.5
$ws.Range("AN8").Value = 6.8

$ws.Range("H9").Value = 2.54
$ws.Range("L9").Value = 1.37
$ws.Range("N9").Value = 3.85
$ws.Range("S9").Value = 2.92
$ws.Range("X9").Value = 18.5
$ws.Range("AA9").Value = 46
$ws.Range("AB9").Value = 13.5
$ws.Range("AC9").Value = 9.800000000000001
$ws.Range("AH9").Value = 19
$ws.Range("AL9").Value = 48
$ws.Range("AO9").Value = 26

$ws.Range("F10").Value = 1.65
$ws.Range("G10").Value = 1.76
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 6.2
$ws.Range("J10").Value = 3.95
$ws.Range("K10").Value = 4.4
$ws.Range("L10").Value = 1.29
$ws.Range("N10").Value = 4.2
$ws.Range("Q10").Value = 1.74
$ws.Range("R10").Value = 1.44
$ws.Range("S10").Value = 2.84
$ws.Range("T10").Value = 1.77
$ws.Range("V10").Value = 1.19
$ws.Range("W10").Value = 2.3
$ws.Range("X10").Value = 21
$ws.Range("Y10").Value = 21
$ws.Range("Z10").Value = 1000
$ws.Range("AB10").Value = 12.5
$ws.Range("AC10").Value = 11.5
$ws.Range("AD10").Value = 22
$ws.Range("AE10").Value = 1000
$ws.Range("AF10").Value = 13.5
$ws.Range("AG10").Value = 13.5
$ws.Range("AH10").Value = 20
$ws.Range("AJ10").Value = 20
$ws.Range("AK10").Value = 18.5
$ws.Range("AL10").Value = 44
$ws.Range("AN10").Value = 11

$ws.Range("L11").Value = 1.25
$ws.Range("X11").Value = 28
$ws.Range("Y11").Value = 34
$ws.Range("AA11").Value = 210
$ws.Range("AB11").Value = 12.5
$ws.Range("AE11").Value = 100
$ws.Range("AG11").Value = 12
$ws.Range("AH11").Value = 24
$ws.Range("AJ11").Value = 16
$ws.Range("AL11").Value = 34
$ws.Range("AN11").Value = 7
$ws.Range("AO11").Value = 100

$ws.Range("F12").Value = 2.12
$ws.Range("G12").Value = 2.18
$ws.Range("H12").Value = 3.45
$ws.Range("I12").Value = 3.65
$ws.Range("L12").Value = 1.26
$ws.Range("R12").Value = 1.58
$ws.Range("V12").Value = 1.37
$ws.Range("X12").Value = 26
$ws.Range("AD12").Value = 18
$ws.Range("AF12").Value = 20
$ws.Range("AH12").Value = 17.5

$ws.Range("F13").Value = 12
$ws.Range("K13").Value = 8.4
$ws.Range("R13").Value = 2.22
$ws.Range("S13").Value = 1.76
$ws.Range("X13").Value = 55
$ws.Range("Y13").Value = 17.5
$ws.Range("AI13").Value = 24

$ws.Range("H14").Value = 4.2
$ws.Range("L14").Value = 1.31
$ws.Range("N14").Value = 5.3
$ws.Range("O14").Value = 1.21
$ws.Range("P14").Value = 2.42
$ws.Range("Q14").Value = 1.65
$ws.Range("R14").Value = 1.59
$ws.Range("S14").Value = 2.58
$ws.Range("T14").Value = 1.64
$ws.Range("U14").Value = 2.5
$ws.Range("X14").Value = 21
$ws.Range("Y14").Value = 21
$ws.Range("Z14").Value = 32
$ws.Range("AA14").Value = 85
$ws.Range("AB14").Value = 12.5
$ws.Range("AC14").Value = 9.199999999999999
$ws.Range("AD14").Value = 17.5
$ws.Range("AF14").Value = 13.5
$ws.Range("AG14").Value = 11
$ws.Range("AH14").Value = 16
$ws.Range("AK14").Value = 17.5
$ws.Range("AL14").Value = 27
$ws.Range("AN14").Value = 9
$ws.Range("AO14").Value = 34

$ws.Range("F15").Value = 1.8
$ws.Range("G15").Value = 1.81
$ws.Range("H15").Value = 4.8
$ws.Range("I15").Value = 4.9
$ws.Range("V15").Value = 1.25
$ws.Range("W15").Value = 2.22

$ws.Range("F16").Value = 1.16
$ws.Range("G16").Value = 1.17
$ws.Range("H16").Value = 26
$ws.Range("I16").Value = 27
$ws.Range("L16").Value = 1.27
$ws.Range("N16").Value = 5.7
$ws.Range("O16").Value = 1.19
$ws.Range("P16").Value = 2.62
$ws.Range("Q16").Value = 1.58
$ws.Range("R16").Value = 1.63
$ws.Range("S16").Value = 2.52
$ws.Range("T16").Value = 2.84
$ws.Range("U16").Value = 1.51
$ws.Range("X16").Value = 29
$ws.Range("Y16").Value = 70
$ws.Range("Z16").Value = 360
$ws.Range("AB16").Value = 9.199999999999999
$ws.Range("AC16").Value = 22
$ws.Range("AF16").Value = 6.8
$ws.Range("AH16").Value = 65
$ws.Range("AI16").Value = 520
$ws.Range("AK16").Value = 16.5
$ws.Range("AN16").Value = 3.55

$ws.Range("N17").Value = 3.9

$ws.Range("H18").Value = 1.49
$ws.Range("Y18").Value = 10.5

$ws.Range("F19").Value = 1.4
$ws.Range("M19").Value = 1.04
$ws.Range("N19").Value = 5.3
$ws.Range("P19").Value = 2.54
$ws.Range("Q19").Value = 1.63
$ws.Range("R19").Value = 1.57
$ws.Range("S19").Value = 2.62
$ws.Range("T19").Value = 1.91
$ws.Range("X19").Value = 27
$ws.Range("Y19").Value = 36
$ws.Range("AA19").Value = 320
$ws.Range("AB19").Value = 11
$ws.Range("AC19").Value = 13.5
$ws.Range("AD19").Value = 36
$ws.Range("AF19").Value = 10.5
$ws.Range("AG19").Value = 12
$ws.Range("AH19").Value = 30
$ws.Range("AI19").Value = 120
$ws.Range("AJ19").Value = 13.5
$ws.Range("AK19").Value = 17.5
$ws.Range("AL19").Value = 40
$ws.Range("AN19").Value = 6.4

$ws.Range("K20").Value = 3.8
$ws.Range("L20").Value = 1.34
$ws.Range("N20").Value = 3.85
$ws.Range("P20").Value = 1.99
$ws.Range("R20").Value = 1.38
$ws.Range("U20").Value = 2.22
$ws.Range("X20").Value = 17
$ws.Range("Y20").Value = 13.5
$ws.Range("Z20").Value = 22
$ws.Range("AA20").Value = 44
$ws.Range("AB20").Value = 13.5
$ws.Range("AE20").Value = 36
$ws.Range("AF20").Value = 22
$ws.Range("AH20").Value = 17.5
$ws.Range("AI20").Value = 46
$ws.Range("AJ20").Value = 48
$ws.Range("AK20").Value = 36
$ws.Range("AL20").Value = 48
$ws.Range("AM20").Value = 95
$ws.Range("AN20").Value = 29
$ws.Range("AO20").Value = 32

$ws.Range("I21").Value = 9.4
$ws.Range("V21").Value = 1.12
$ws.Range("AL21").Value = 40
